# Refresh cached Universalis market-price figures for the Leviathan data
# center across each Disciple of the Hand sheet (scheduled runner update).
# Values below are the latest pull; profit columns are recomputed from them.

$wb = $excel.ActiveWorkbook

$updates = @(
    @{ Sheet = "ALC"; Cell = "H40"; Value = 3759.1667 }
    @{ Sheet = "ALC"; Cell = "I40"; Value = 1961.3334 }
    @{ Sheet = "ALC"; Cell = "J40"; Value = 4658.0835 }
    @{ Sheet = "ALC"; Cell = "K40"; Value = 1961.3334 }
    @{ Sheet = "ALC"; Cell = "L40"; Value = 4658.0835 }
    @{ Sheet = "ALC"; Cell = "M40"; Value = -1786.3334 }
    @{ Sheet = "ALC"; Cell = "N40"; Value = -5008.0835 }
    @{ Sheet = "ALC"; Cell = "H41"; Value = 433.91666 }
    @{ Sheet = "ALC"; Cell = "I41"; Value = 461.9 }
    @{ Sheet = "ALC"; Cell = "J41"; Value = 413.92856 }
    @{ Sheet = "ALC"; Cell = "K41"; Value = 461.9 }
    @{ Sheet = "ALC"; Cell = "L41"; Value = 413.92856 }
    @{ Sheet = "ALC"; Cell = "M41"; Value = -21.89999999999998 }
    @{ Sheet = "ALC"; Cell = "N41"; Value = -1293.92856 }
    @{ Sheet = "ALC"; Cell = "H70"; Value = 13550.375 }
    @{ Sheet = "ALC"; Cell = "I70"; Value = 15314.286 }
    @{ Sheet = "ALC"; Cell = "K70"; Value = 45942.858 }
    @{ Sheet = "ALC"; Cell = "M70"; Value = -45672.858 }
    @{ Sheet = "ALC"; Cell = "H73"; Value = 13550.375 }
    @{ Sheet = "ALC"; Cell = "I73"; Value = 15314.286 }
    @{ Sheet = "ALC"; Cell = "K73"; Value = 45942.858 }
    @{ Sheet = "ALC"; Cell = "M73"; Value = -45006.858 }
    @{ Sheet = "ALC"; Cell = "H112"; Value = 1898.0741 }
    @{ Sheet = "ALC"; Cell = "J112"; Value = 1673.92 }
    @{ Sheet = "ALC"; Cell = "L112"; Value = 5021.76 }
    @{ Sheet = "ALC"; Cell = "N112"; Value = -7237.76 }
    @{ Sheet = "ALC"; Cell = "H132"; Value = 3548.2068 }
    @{ Sheet = "ALC"; Cell = "I132"; Value = 1575.6522 }
    @{ Sheet = "ALC"; Cell = "K132"; Value = 4726.9566 }
    @{ Sheet = "ALC"; Cell = "M132"; Value = -2196.9566 }
    @{ Sheet = "ALC"; Cell = "H137"; Value = 4232.175 }
    @{ Sheet = "ALC"; Cell = "I137"; Value = 4358.913 }
    @{ Sheet = "ALC"; Cell = "J137"; Value = 4060.7058 }
    @{ Sheet = "ALC"; Cell = "K137"; Value = 13076.739 }
    @{ Sheet = "ALC"; Cell = "L137"; Value = 12182.1174 }
    @{ Sheet = "ALC"; Cell = "M137"; Value = -10526.739 }
    @{ Sheet = "ALC"; Cell = "N137"; Value = -17282.1174 }
    @{ Sheet = "ALC"; Cell = "H138"; Value = 3112.068 }
    @{ Sheet = "ALC"; Cell = "I138"; Value = 2619.7 }
    @{ Sheet = "ALC"; Cell = "J138"; Value = 3256.8823 }
    @{ Sheet = "ALC"; Cell = "K138"; Value = 7859.099999999999 }
    @{ Sheet = "ALC"; Cell = "L138"; Value = 9770.6469 }
    @{ Sheet = "ALC"; Cell = "M138"; Value = -2719.099999999999 }
    @{ Sheet = "ALC"; Cell = "N138"; Value = -20050.6469 }
    @{ Sheet = "ARM"; Cell = "H32"; Value = 5441568 }
    @{ Sheet = "ARM"; Cell = "I32"; Value = 750948 }
    @{ Sheet = "ARM"; Cell = "J32"; Value = 41793870 }
    @{ Sheet = "ARM"; Cell = "K32"; Value = 750948 }
    @{ Sheet = "ARM"; Cell = "L32"; Value = 41793870 }
    @{ Sheet = "ARM"; Cell = "M32"; Value = -750661 }
    @{ Sheet = "ARM"; Cell = "N32"; Value = -41794444 }
    @{ Sheet = "ARM"; Cell = "H61"; Value = 1186.3636 }
    @{ Sheet = "ARM"; Cell = "I61"; Value = 1215 }
    @{ Sheet = "ARM"; Cell = "J61"; Value = 900 }
    @{ Sheet = "ARM"; Cell = "K61"; Value = 1215 }
    @{ Sheet = "ARM"; Cell = "L61"; Value = 900 }
    @{ Sheet = "ARM"; Cell = "M61"; Value = -1003 }
    @{ Sheet = "ARM"; Cell = "N61"; Value = -1324 }
    @{ Sheet = "ARM"; Cell = "H74"; Value = 2169.4285 }
    @{ Sheet = "ARM"; Cell = "I74"; Value = 1577.8667 }
    @{ Sheet = "ARM"; Cell = "J74"; Value = 3648.3333 }
    @{ Sheet = "ARM"; Cell = "K74"; Value = 1577.8667 }
    @{ Sheet = "ARM"; Cell = "L74"; Value = 3648.3333 }
    @{ Sheet = "ARM"; Cell = "M74"; Value = -703.8667 }
    @{ Sheet = "ARM"; Cell = "N74"; Value = -5396.3333 }
    @{ Sheet = "ARM"; Cell = "H77"; Value = 2169.4285 }
    @{ Sheet = "ARM"; Cell = "I77"; Value = 1577.8667 }
    @{ Sheet = "ARM"; Cell = "J77"; Value = 3648.3333 }
    @{ Sheet = "ARM"; Cell = "K77"; Value = 7889.333500000001 }
    @{ Sheet = "ARM"; Cell = "L77"; Value = 18241.6665 }
    @{ Sheet = "ARM"; Cell = "M77"; Value = -3521.333500000001 }
    @{ Sheet = "ARM"; Cell = "N77"; Value = -26977.6665 }
    @{ Sheet = "ARM"; Cell = "H102"; Value = 2356.3333 }
    @{ Sheet = "ARM"; Cell = "I102"; Value = 2260.0833 }
    @{ Sheet = "ARM"; Cell = "J102"; Value = 2741.3333 }
    @{ Sheet = "ARM"; Cell = "K102"; Value = 2260.0833 }
    @{ Sheet = "ARM"; Cell = "L102"; Value = 2741.3333 }
    @{ Sheet = "ARM"; Cell = "M102"; Value = -638.0832999999998 }
    @{ Sheet = "ARM"; Cell = "N102"; Value = -5985.3333 }
    @{ Sheet = "ARM"; Cell = "H110"; Value = 1122.091 }
    @{ Sheet = "ARM"; Cell = "I110"; Value = 937.4286 }
    @{ Sheet = "ARM"; Cell = "K110"; Value = 937.4286 }
    @{ Sheet = "ARM"; Cell = "M110"; Value = 1107.5714 }
    @{ Sheet = "ARM"; Cell = "H119"; Value = 0 }
    @{ Sheet = "ARM"; Cell = "J119"; Value = 0 }
    @{ Sheet = "ARM"; Cell = "L119"; Value = 0 }
    @{ Sheet = "ARM"; Cell = "N119"; Value = $null }
    @{ Sheet = "ARM"; Cell = "H136"; Value = 1186.3636 }
    @{ Sheet = "ARM"; Cell = "I136"; Value = 1215 }
    @{ Sheet = "ARM"; Cell = "J136"; Value = 900 }
    @{ Sheet = "ARM"; Cell = "K136"; Value = 3645 }
    @{ Sheet = "ARM"; Cell = "L136"; Value = 2700 }
    @{ Sheet = "ARM"; Cell = "M136"; Value = -1095 }
    @{ Sheet = "ARM"; Cell = "N136"; Value = -7800 }
    @{ Sheet = "BSM"; Cell = "H86"; Value = 2729.6 }
    @{ Sheet = "BSM"; Cell = "J86"; Value = 3046 }
    @{ Sheet = "BSM"; Cell = "L86"; Value = 3046 }
    @{ Sheet = "BSM"; Cell = "N86"; Value = -5292 }
    @{ Sheet = "BSM"; Cell = "H89"; Value = 2729.6 }
    @{ Sheet = "BSM"; Cell = "J89"; Value = 3046 }
    @{ Sheet = "BSM"; Cell = "L89"; Value = 15230 }
    @{ Sheet = "BSM"; Cell = "N89"; Value = -26462 }
    @{ Sheet = "BSM"; Cell = "H94"; Value = 34709.61 }
    @{ Sheet = "BSM"; Cell = "I94"; Value = 1642.45 }
    @{ Sheet = "BSM"; Cell = "K94"; Value = 1642.45 }
    @{ Sheet = "BSM"; Cell = "M94"; Value = -1191.45 }
    @{ Sheet = "BSM"; Cell = "H105"; Value = 49307.445 }
    @{ Sheet = "BSM"; Cell = "I105"; Value = 3750 }
    @{ Sheet = "BSM"; Cell = "K105"; Value = 3750 }
    @{ Sheet = "BSM"; Cell = "M105"; Value = -2003 }
    @{ Sheet = "BSM"; Cell = "H141"; Value = 99977 }
    @{ Sheet = "BSM"; Cell = "J141"; Value = 99977 }
    @{ Sheet = "BSM"; Cell = "L141"; Value = 99977 }
    @{ Sheet = "BSM"; Cell = "N141"; Value = -110337 }
    @{ Sheet = "CRP"; Cell = "H26"; Value = 8000 }
    @{ Sheet = "CRP"; Cell = "J26"; Value = 8000 }
    @{ Sheet = "CRP"; Cell = "L26"; Value = 8000 }
    @{ Sheet = "CRP"; Cell = "N26"; Value = -8574 }
    @{ Sheet = "CRP"; Cell = "H31"; Value = 4065.963 }
    @{ Sheet = "CRP"; Cell = "I31"; Value = 2379.5715 }
    @{ Sheet = "CRP"; Cell = "J31"; Value = 5882.077 }
    @{ Sheet = "CRP"; Cell = "K31"; Value = 2379.5715 }
    @{ Sheet = "CRP"; Cell = "L31"; Value = 5882.077 }
    @{ Sheet = "CRP"; Cell = "M31"; Value = -2084.5715 }
    @{ Sheet = "CRP"; Cell = "N31"; Value = -6472.077 }
    @{ Sheet = "CRP"; Cell = "H34"; Value = 4065.963 }
    @{ Sheet = "CRP"; Cell = "I34"; Value = 2379.5715 }
    @{ Sheet = "CRP"; Cell = "J34"; Value = 5882.077 }
    @{ Sheet = "CRP"; Cell = "K34"; Value = 2379.5715 }
    @{ Sheet = "CRP"; Cell = "L34"; Value = 5882.077 }
    @{ Sheet = "CRP"; Cell = "M34"; Value = -2177.5715 }
    @{ Sheet = "CRP"; Cell = "N34"; Value = -6286.077 }
    @{ Sheet = "CRP"; Cell = "H51"; Value = 12499 }
    @{ Sheet = "CRP"; Cell = "J51"; Value = 13180.728 }
    @{ Sheet = "CRP"; Cell = "L51"; Value = 13180.728 }
    @{ Sheet = "CRP"; Cell = "N51"; Value = -14652.728 }
    @{ Sheet = "CRP"; Cell = "H59"; Value = 20000 }
    @{ Sheet = "CRP"; Cell = "J59"; Value = 20000 }
    @{ Sheet = "CRP"; Cell = "L59"; Value = 20000 }
    @{ Sheet = "CRP"; Cell = "N59"; Value = -22290 }
    @{ Sheet = "CRP"; Cell = "H61"; Value = 12499 }
    @{ Sheet = "CRP"; Cell = "J61"; Value = 13180.728 }
    @{ Sheet = "CRP"; Cell = "L61"; Value = 13180.728 }
    @{ Sheet = "CRP"; Cell = "N61"; Value = -13876.728 }
    @{ Sheet = "CRP"; Cell = "H97"; Value = 28678.8 }
    @{ Sheet = "CRP"; Cell = "I97"; Value = 25000 }
    @{ Sheet = "CRP"; Cell = "J97"; Value = 29598.5 }
    @{ Sheet = "CRP"; Cell = "K97"; Value = 25000 }
    @{ Sheet = "CRP"; Cell = "L97"; Value = 29598.5 }
    @{ Sheet = "CRP"; Cell = "M97"; Value = -24009 }
    @{ Sheet = "CRP"; Cell = "N97"; Value = -31580.5 }
    @{ Sheet = "CRP"; Cell = "H99"; Value = 31107 }
    @{ Sheet = "CRP"; Cell = "I99"; Value = 34545.5 }
    @{ Sheet = "CRP"; Cell = "K99"; Value = 34545.5 }
    @{ Sheet = "CRP"; Cell = "M99"; Value = -33047.5 }
    @{ Sheet = "CRP"; Cell = "H107"; Value = 1633.8889 }
    @{ Sheet = "CRP"; Cell = "J107"; Value = 1972.375 }
    @{ Sheet = "CRP"; Cell = "L107"; Value = 1972.375 }
    @{ Sheet = "CRP"; Cell = "N107"; Value = -5812.375 }
    @{ Sheet = "CRP"; Cell = "H114"; Value = 40000 }
    @{ Sheet = "CRP"; Cell = "J114"; Value = 40000 }
    @{ Sheet = "CRP"; Cell = "L114"; Value = 40000 }
    @{ Sheet = "CRP"; Cell = "N114"; Value = -48678 }
    @{ Sheet = "CRP"; Cell = "H126"; Value = 31107 }
    @{ Sheet = "CRP"; Cell = "I126"; Value = 34545.5 }
    @{ Sheet = "CRP"; Cell = "K126"; Value = 103636.5 }
    @{ Sheet = "CRP"; Cell = "M126"; Value = -101166.5 }
    @{ Sheet = "CRP"; Cell = "H132"; Value = 2817.9167 }
    @{ Sheet = "CRP"; Cell = "I132"; Value = 2694.9333 }
    @{ Sheet = "CRP"; Cell = "K132"; Value = 8084.7999 }
    @{ Sheet = "CRP"; Cell = "M132"; Value = -5554.7999 }
    @{ Sheet = "CRP"; Cell = "H133"; Value = 0 }
    @{ Sheet = "CRP"; Cell = "J133"; Value = 0 }
    @{ Sheet = "CRP"; Cell = "L133"; Value = 0 }
    @{ Sheet = "CRP"; Cell = "N133"; Value = $null }
    @{ Sheet = "CRP"; Cell = "H134"; Value = 5296.5386 }
    @{ Sheet = "CRP"; Cell = "I134"; Value = 4777.909 }
    @{ Sheet = "CRP"; Cell = "K134"; Value = 14333.727 }
    @{ Sheet = "CRP"; Cell = "M134"; Value = -11798.727 }
    @{ Sheet = "CUL"; Cell = "H97"; Value = 394.72726 }
    @{ Sheet = "CUL"; Cell = "I97"; Value = 431.5 }
    @{ Sheet = "CUL"; Cell = "K97"; Value = 1294.5 }
    @{ Sheet = "CUL"; Cell = "M97"; Value = -798.5 }
    @{ Sheet = "CUL"; Cell = "H107"; Value = 878.41174 }
    @{ Sheet = "CUL"; Cell = "I107"; Value = 448.22223 }
    @{ Sheet = "CUL"; Cell = "J107"; Value = 1362.375 }
    @{ Sheet = "CUL"; Cell = "K107"; Value = 1344.66669 }
    @{ Sheet = "CUL"; Cell = "L107"; Value = 4087.125 }
    @{ Sheet = "CUL"; Cell = "M107"; Value = 575.33331 }
    @{ Sheet = "CUL"; Cell = "N107"; Value = -7927.125 }
    @{ Sheet = "CUL"; Cell = "H137"; Value = 7147232.5 }
    @{ Sheet = "CUL"; Cell = "J137"; Value = 6141.375 }
    @{ Sheet = "CUL"; Cell = "L137"; Value = 18424.125 }
    @{ Sheet = "CUL"; Cell = "N137"; Value = -28624.125 }
    @{ Sheet = "GSM"; Cell = "H70"; Value = 41231748 }
    @{ Sheet = "GSM"; Cell = "I70"; Value = 5499.7 }
    @{ Sheet = "GSM"; Cell = "K70"; Value = 5499.7 }
    @{ Sheet = "GSM"; Cell = "M70"; Value = -5229.7 }
    @{ Sheet = "GSM"; Cell = "H73"; Value = 41231748 }
    @{ Sheet = "GSM"; Cell = "I73"; Value = 5499.7 }
    @{ Sheet = "GSM"; Cell = "K73"; Value = 5499.7 }
    @{ Sheet = "GSM"; Cell = "M73"; Value = -4563.7 }
    @{ Sheet = "GSM"; Cell = "H132"; Value = 4152.593 }
    @{ Sheet = "GSM"; Cell = "I132"; Value = 4767.222 }
    @{ Sheet = "GSM"; Cell = "K132"; Value = 14301.666 }
    @{ Sheet = "GSM"; Cell = "M132"; Value = -11771.666 }
    @{ Sheet = "WVR"; Cell = "H62"; Value = 128999.43 }
    @{ Sheet = "WVR"; Cell = "I62"; Value = 2750.25 }
    @{ Sheet = "WVR"; Cell = "J62"; Value = 179499.1 }
    @{ Sheet = "WVR"; Cell = "K62"; Value = 2750.25 }
    @{ Sheet = "WVR"; Cell = "L62"; Value = 179499.1 }
    @{ Sheet = "WVR"; Cell = "M62"; Value = -2126.25 }
    @{ Sheet = "WVR"; Cell = "N62"; Value = -180747.1 }
    @{ Sheet = "WVR"; Cell = "H65"; Value = 128999.43 }
    @{ Sheet = "WVR"; Cell = "I65"; Value = 2750.25 }
    @{ Sheet = "WVR"; Cell = "J65"; Value = 179499.1 }
    @{ Sheet = "WVR"; Cell = "K65"; Value = 13751.25 }
    @{ Sheet = "WVR"; Cell = "L65"; Value = 897495.5 }
    @{ Sheet = "WVR"; Cell = "M65"; Value = -10631.25 }
    @{ Sheet = "WVR"; Cell = "N65"; Value = -903735.5 }
)

foreach ($u in $updates) {
    $ws = $wb.Worksheets.Item($u.Sheet)
    $rng = $ws.Range($u.Cell)
    if ($null -eq $u.Value) {
        $rng.ClearContents()
    } else {
        $rng.Value = $u.Value
    }
}
